$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column Z (26th column), shifting
# existing columns Z:AJ to AB:AL. This makes room for the new
# "derived_variable" / "derivation_description" header columns.
$ws.Range("Z1:AA1").EntireColumn.Insert()

# New header cells (row 1) - headers use the bold/bordered style already
# present on the row, so copy formatting from the neighbouring header cell.
$ws.Range("Y1").Copy()
$ws.Range("Z1:AA1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("Z1").Value = "derived_variable"
$ws.Range("AA1").Value = "derivation_description"

# Update view state to match the authored selection/scroll position.
$ws.Application.ActiveWindow.ScrollColumn = 16
$ws.Range("Y14").Select()
